$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value, exactly as in the target diff.
# Column D holds price text that often LOOKS numeric ("586.93", "1.00", "0.0360", ...).
# Force those ranges to Text format first so Excel keeps the literal string
# (no float coercion / lost trailing zeros), matching the inlineStr cells in the workbook.
$updates = [ordered]@{
    "D2" = "69.625.55"
    "E2" = "  +5.07%  "
    "D3" = "3.433.99"
    "E3" = "  +11.45%  "
    "E4" = "  +0.03%  "
    "D5" = "586.93"
    "E5" = "  +2.09%  "
    "D6" = "184.46"
    "E6" = "  +8.68%  "
    "E7" = "  +0.01%  "
    "D8" = "3.424.10"
    "E8" = "  +11.23%  "
    "D9" = "0.530"
    "E9" = "  +4.15%  "
    "E10" = "  +3.48%  "
    "E11" = "  +4.68%  "
    "D12" = "0.484"
    "E12" = "  +2.95%  "
    "B13" = "ShibaInu"
    "C13" = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
    "D13" = "0.0000248"
    "E13" = "  +3.82%  "
    "B14" = "Avalanche"
    "C14" = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
    "D14" = "38.19"
    "E14" = "  +6.77%  "
    "D15" = "4.006.50"
    "E15" = "  +11.51%  "
    "D16" = "69.724.04"
    "E16" = "  +5.30%  "
    "E17" = "  +1.14%  "
    "D18" = "3.433.33"
    "E18" = "  +11.53%  "
    "D19" = "7.34"
    "E19" = "  +5.46%  "
    "D20" = "16.83"
    "E20" = "  +0.72%  "
    "D21" = "500.13"
    "E21" = "  +2.64%  "
    "D22" = "8.67"
    "E22" = "  +12.47%  "
    "D23" = "0.722"
    "E23" = "  +5.13%  "
    "D24" = "86.28"
    "E24" = "  +4.53%  "
    "E25" = "  +4.21%  "
    "D26" = "2.35"
    "E26" = "  +6.87%  "
    "D27" = "10.70"
    "E27" = "  +4.56%  "
    "D28" = "1.00"
    "E28" = "  -0.05%  "
    "B29" = "ImmutableX"
    "C29" = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
    "D29" = "2.50"
    "E29" = "  +10.77%  "
    "B30" = "NEARProtocol"
    "C30" = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
    "D30" = "8.14"
    "E30" = "  +3.42%  "
    "E31" = "  +3.85%  "
    "D32" = "30.11"
    "E32" = "  +8.56%  "
    "E33" = "  +15.46%  "
    "D34" = "0.115"
    "E34" = "  +4.00%  "
    "E35" = "  +0.09%  "
    "D36" = "6.05"
    "E36" = "  +8.57%  "
    "E37" = "  +6.18%  "
    "D38" = "47.75"
    "E38" = "  +1.07%  "
    "E39" = "  +9.29%  "
    "E40" = "  +6.55%  "
    "D41" = "0.128"
    "E41" = "  +4.75%  "
    "D42" = "50.13"
    "E42" = "  +2.14%  "
    "D43" = "8.63"
    "E43" = "  +4.32%  "
    "D44" = "414.03"
    "E44" = "  +13.43%  "
    "D45" = "2.81"
    "E45" = "  +12.43%  "
    "D46" = "2.932.85"
    "E46" = "  +5.30%  "
    "D47" = "27.93"
    "E47" = "  +14.39%  "
    "D48" = "0.0360"
    "E48" = "  +4.58%  "
    "D49" = "134.37"
    "E49" = "  -0.15%  "
    "E50" = "  +0.04%  "
    "D51" = "2.43"
    "E51" = "  +12.87%  "
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    if ($ref.StartsWith("D")) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $updates[$ref]
}
